# Applies Week 15 simulation updates to Raiders Players Data workbook.
$wb = $excel.ActiveWorkbook

# --- Rushing sheet ---
$rushing = $wb.Worksheets.Item("Rushing")
$rushing.Range("C5").Value = 66
$rushing.Range("D5").Value = 51

# --- Receiving sheet ---
$receiving = $wb.Worksheets.Item("Receiving")
$receiving.Range("C2").Value = 47
$receiving.Range("D2").Value = 41

$receiving.Range("E7").Value = 23

$receiving.Range("C8").Value = 92
$receiving.Range("D8").Value = 76
$receiving.Range("G8").Value = 15
$receiving.Range("H8").Value = 11

$receiving.Range("C9").Value = 16
$receiving.Range("D9").Value = 12
$receiving.Range("E9").Value = 15

$receiving.Range("C11").Value = 4
$receiving.Range("D11").Value = 3
$receiving.Range("E11").Value = 4

$receiving.Range("C13").Value = 22
$receiving.Range("D13").Value = 13
$receiving.Range("G13").Value = 6
